$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = '北方稀土'
$ws.Range('B2').Value = '北方稀土'
$ws.Range('C2').Value = '北方稀土'
$ws.Range('A3').Value = '凯美特气'
$ws.Range('B3').Value = '包钢股份'
$ws.Range('C3').Value = '海南华铁'
$ws.Range('A4').Value = '合锻智能'
$ws.Range('B4').Value = '长城军工'
$ws.Range('C4').Value = '闻泰科技'
$ws.Range('A5').Value = '至纯科技'
$ws.Range('B5').Value = '闻泰科技'
$ws.Range('C5').Value = '赣锋锂业'
$ws.Range('A6').Value = '包钢股份'
$ws.Range('B6').Value = '至纯科技'
$ws.Range('C6').Value = '永鼎股份'
$ws.Range('A7').Value = '新莱应材'
$ws.Range('B7').Value = '凯美特气'
$ws.Range('C7').Value = '白银有色'
$ws.Range('A8').Value = '中国稀土'
$ws.Range('B8').Value = '合锻智能'
$ws.Range('C8').Value = '通富微电'
$ws.Range('A9').Value = '安泰科技'
$ws.Range('B9').Value = '安泰科技'
$ws.Range('C9').Value = '山子高科'
$ws.Range('A10').Value = '闻泰科技'
$ws.Range('B10').Value = '新莱应材'
$ws.Range('C10').Value = '楚江新材'
$ws.Range('A11').Value = '长城军工'
$ws.Range('B11').Value = '永鼎股份'
$ws.Range('C11').Value = '凯美特气'
$ws.Range('A12').Value = '永鼎股份'
$ws.Range('B12').Value = '中国稀土'
$ws.Range('C12').Value = '合锻智能'
$ws.Range('A13').Value = '赣锋锂业'
$ws.Range('B13').Value = '赣锋锂业'
$ws.Range('C13').Value = '长城军工'
$ws.Range('A14').Value = '通富微电'
$ws.Range('B14').Value = '通富微电'
$ws.Range('C14').Value = '至纯科技'
$ws.Range('A15').Value = '天际股份'
$ws.Range('B15').Value = '上海电气'
$ws.Range('C15').Value = '中国稀土'
$ws.Range('A16').Value = '海南华铁'
$ws.Range('B16').Value = '银河磁体'
$ws.Range('C16').Value = '安泰科技'
$ws.Range('A17').Value = '山子高科'
$ws.Range('B17').Value = '金力永磁'
$ws.Range('C17').Value = '包钢股份'
$ws.Range('A18').Value = '中国软件'
$ws.Range('B18').Value = '山子高科'
$ws.Range('C18').Value = '上海电气'
$ws.Range('A19').Value = '金力永磁'
$ws.Range('B19').Value = '中国软件'
$ws.Range('C19').Value = '新莱应材'
$ws.Range('A20').Value = '多氟多'
$ws.Range('B20').Value = '海南华铁'
$ws.Range('C20').Value = '金力永磁'
$ws.Range('A21').Value = '白银有色'
$ws.Range('B21').Value = '白银有色'
$ws.Range('C21').Value = '华友钴业'
